# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect the latest generated output (scraped data refresh).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 92
    4  = 1534
    5  = 586
    6  = 1081
    7  = 11200
    8  = 5
    10 = 36
    11 = 331
    12 = 1075
    13 = 770
    14 = 12273
    15 = 12873
    16 = 32
    17 = 130
    22 = 63
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
